$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 8741.666999999999
$ws.Range("I40").Value = 8210.263000000001
$ws.Range("J40").Value = 9659.546
$ws.Range("K40").Value = 8210.263000000001
$ws.Range("L40").Value = 9659.546
$ws.Range("M40").Value = -8035.263000000001
$ws.Range("N40").Value = -10009.546
# Row 70
$ws.Range("H70").Value = 6414.1
$ws.Range("I70").Value = 5344.222
$ws.Range("J70").Value = 7289.4546
$ws.Range("K70").Value = 16032.666
$ws.Range("L70").Value = 21868.3638
$ws.Range("M70").Value = -15762.666
$ws.Range("N70").Value = -22408.3638
# Row 73
$ws.Range("H73").Value = 6414.1
$ws.Range("I73").Value = 5344.222
$ws.Range("J73").Value = 7289.4546
$ws.Range("K73").Value = 16032.666
$ws.Range("L73").Value = 21868.3638
$ws.Range("M73").Value = -15096.666
$ws.Range("N73").Value = -23740.3638
# Row 76
$ws.Range("H76").Value = 2000
$ws.Range("I76").Value = 2000
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 2000
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -1685
$ws.Range("N76").ClearContents()
# Row 79
$ws.Range("H79").Value = 2000
$ws.Range("I79").Value = 2000
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 2000
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -908
$ws.Range("N79").ClearContents()
# Row 86
$ws.Range("H86").Value = 3763450.2
$ws.Range("I86").Value = 3246.8572
$ws.Range("K86").Value = 3246.8572
$ws.Range("M86").Value = -2123.8572
# Row 89
$ws.Range("H89").Value = 3763450.2
$ws.Range("I89").Value = 3246.8572
$ws.Range("K89").Value = 16234.286
$ws.Range("M89").Value = -10618.286
# Row 113
$ws.Range("H113").Value = 11974.25
$ws.Range("I113").Value = 7438.4
$ws.Range("J113").Value = 15214.143
$ws.Range("K113").Value = 7438.4
$ws.Range("L113").Value = 15214.143
$ws.Range("M113").Value = -4184.4
$ws.Range("N113").Value = -21722.143
# Row 138
$ws.Range("H138").Value = 4251.161
$ws.Range("I138").Value = 3408.6
$ws.Range("J138").Value = 4652.381
$ws.Range("K138").Value = 10225.8
$ws.Range("L138").Value = 13957.143
$ws.Range("M138").Value = -5085.799999999999
$ws.Range("N138").Value = -24237.143
# Row 141
$ws.Range("H141").Value = 6739.364
$ws.Range("I141").Value = 4407.107
$ws.Range("K141").Value = 13221.321
$ws.Range("M141").Value = -8041.321

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 14498907
$ws.Range("I32").Value = 15157824
$ws.Range("K32").Value = 15157824
$ws.Range("M32").Value = -15157537
# Row 61
$ws.Range("H61").Value = 4520.4287
$ws.Range("I61").Value = 3607.1667
$ws.Range("K61").Value = 3607.1667
$ws.Range("M61").Value = -3395.1667
# Row 74
$ws.Range("H74").Value = 1507.1904
$ws.Range("I74").Value = 1322.5385
$ws.Range("J74").Value = 1807.25
$ws.Range("K74").Value = 1322.5385
$ws.Range("L74").Value = 1807.25
$ws.Range("M74").Value = -448.5385000000001
$ws.Range("N74").Value = -3555.25
# Row 77
$ws.Range("H77").Value = 1507.1904
$ws.Range("I77").Value = 1322.5385
$ws.Range("J77").Value = 1807.25
$ws.Range("K77").Value = 6612.692500000001
$ws.Range("L77").Value = 9036.25
$ws.Range("M77").Value = -2244.692500000001
$ws.Range("N77").Value = -17772.25
# Row 122
$ws.Range("H122").Value = 6095
$ws.Range("I122").Value = 6333.5713
$ws.Range("K122").Value = 19000.7139
$ws.Range("M122").Value = -16550.7139
# Row 132
$ws.Range("H132").Value = 6834
$ws.Range("I132").Value = 6401.727
$ws.Range("K132").Value = 19205.181
$ws.Range("M132").Value = -16675.181
# Row 136
$ws.Range("H136").Value = 4520.4287
$ws.Range("I136").Value = 3607.1667
$ws.Range("K136").Value = 10821.5001
$ws.Range("M136").Value = -8271.500100000001

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 1045
$ws.Range("I22").Value = 856.25
$ws.Range("K22").Value = 856.25
$ws.Range("M22").Value = -683.25

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 7508.65
$ws.Range("I16").Value = 5545.3335
$ws.Range("K16").Value = 5545.3335
$ws.Range("M16").Value = -5258.3335
# Row 48
$ws.Range("H48").Value = 42495
$ws.Range("J48").Value = 42495
$ws.Range("L48").Value = 42495
$ws.Range("N48").Value = -43447
# Row 109
$ws.Range("H109").Value = 27714.143
$ws.Range("J109").Value = 27714.143
$ws.Range("L109").Value = 27714.143
$ws.Range("N109").Value = -29794.143
# Row 113
$ws.Range("H113").Value = 7508.65
$ws.Range("I113").Value = 5545.3335
$ws.Range("K113").Value = 5545.3335
$ws.Range("M113").Value = -3375.3335
# Row 132
$ws.Range("H132").Value = 978
$ws.Range("I132").Value = 978
$ws.Range("K132").Value = 2934
$ws.Range("M132").Value = -404

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1402.1
$ws.Range("I5").Value = 1553.3846
$ws.Range("J5").Value = 1121.1428
$ws.Range("K5").Value = 4660.1538
$ws.Range("L5").Value = 3363.4284
$ws.Range("M5").Value = -4548.1538
$ws.Range("N5").Value = -3587.4284
# Row 68
$ws.Range("H68").Value = 1833
$ws.Range("I68").Value = 1399.8
$ws.Range("J68").Value = 2374.5
$ws.Range("K68").Value = 4199.4
$ws.Range("L68").Value = 7123.5
$ws.Range("M68").Value = -3388.4
$ws.Range("N68").Value = -8745.5
# Row 71
$ws.Range("H71").Value = 1833
$ws.Range("I71").Value = 1399.8
$ws.Range("J71").Value = 2374.5
$ws.Range("K71").Value = 12598.2
$ws.Range("L71").Value = 21370.5
$ws.Range("M71").Value = -8542.199999999999
$ws.Range("N71").Value = -29482.5
# Row 107
$ws.Range("H107").Value = 2944.3809
$ws.Range("I107").Value = 2890.6365
$ws.Range("J107").Value = 3003.5
$ws.Range("K107").Value = 8671.9095
$ws.Range("L107").Value = 9010.5
$ws.Range("M107").Value = -6751.9095
$ws.Range("N107").Value = -12850.5
# Row 131
$ws.Range("H131").Value = 1594.1428
$ws.Range("I131").Value = 656.6842
$ws.Range("K131").Value = 1970.0526
$ws.Range("M131").Value = 3069.9474
# Row 135
$ws.Range("H135").Value = 1402.1
$ws.Range("I135").Value = 1553.3846
$ws.Range("J135").Value = 1121.1428
$ws.Range("K135").Value = 13980.4614
$ws.Range("L135").Value = 10090.2852
$ws.Range("M135").Value = -11445.4614
$ws.Range("N135").Value = -15160.2852

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 141
$ws.Range("H141").Value = 71000
$ws.Range("J141").Value = 71000
$ws.Range("L141").Value = 71000
$ws.Range("N141").Value = -81360

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 2633.3
$ws.Range("I16").Value = 2614.889
$ws.Range("K16").Value = 2614.889
$ws.Range("M16").Value = -2444.889
# Row 40
$ws.Range("H40").Value = 4997.5
$ws.Range("I40").Value = 4995
$ws.Range("K40").Value = 4995
$ws.Range("M40").Value = -4859
# Row 46
$ws.Range("H46").Value = 1889.72
$ws.Range("J46").Value = 2078.9546
$ws.Range("L46").Value = 2078.9546
$ws.Range("N46").Value = -2454.9546
# Row 68
$ws.Range("H68").Value = 3333
$ws.Range("J68").Value = 3333
$ws.Range("L68").Value = 3333
$ws.Range("N68").Value = -4831
# Row 71
$ws.Range("H71").Value = 3333
$ws.Range("J71").Value = 3333
$ws.Range("L71").Value = 16665
$ws.Range("N71").Value = -24153
# Row 82
$ws.Range("H82").Value = 1486.7333
$ws.Range("I82").Value = 1663.6
$ws.Range("K82").Value = 1663.6
$ws.Range("M82").Value = -1302.6
# Row 85
$ws.Range("H85").Value = 1486.7333
$ws.Range("I85").Value = 1663.6
$ws.Range("K85").Value = 1663.6
$ws.Range("M85").Value = -415.5999999999999
# Row 112
$ws.Range("H112").Value = 29599.9
$ws.Range("J112").Value = 29777.777
$ws.Range("L112").Value = 29777.777
$ws.Range("N112").Value = -32731.777

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 26
$ws.Range("H26").Value = 12512
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()
# Row 62
$ws.Range("H62").Value = 33335952
$ws.Range("I62").Value = 3263
$ws.Range("J62").Value = 83334984
$ws.Range("K62").Value = 3263
$ws.Range("L62").Value = 83334984
$ws.Range("M62").Value = -2639
$ws.Range("N62").Value = -83336232
# Row 65
$ws.Range("H65").Value = 33335952
$ws.Range("I65").Value = 3263
$ws.Range("J65").Value = 83334984
$ws.Range("K65").Value = 16315
$ws.Range("L65").Value = 416674920
$ws.Range("M65").Value = -13195
$ws.Range("N65").Value = -416681160
# Row 109
$ws.Range("H109").Value = 16590.908
$ws.Range("J109").Value = 16590.908
$ws.Range("L109").Value = 16590.908
$ws.Range("N109").Value = -19364.908
# Row 132
$ws.Range("H132").Value = 4059.5217
$ws.Range("I132").Value = 4080.0588
$ws.Range("J132").Value = 4001.3333
$ws.Range("K132").Value = 12240.1764
$ws.Range("L132").Value = 12003.9999
$ws.Range("M132").Value = -9710.1764
$ws.Range("N132").Value = -17063.9999
# Row 136
$ws.Range("H136").Value = 3559.3428
$ws.Range("I136").Value = 2916.8167
$ws.Range("J136").Value = 7414.5
$ws.Range("K136").Value = 8750.4501
$ws.Range("L136").Value = 22243.5
$ws.Range("M136").Value = -6200.4501
$ws.Range("N136").Value = -27343.5
